$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: new date/experiment-note header, same style as A32 ---
$ws.Range("A42").Value = "2024/08/18（直接フォースゲージつけた）"
$ws.Range("A32").Copy()
$ws.Range("A42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 43: labels / slope coefficients (copied from row 33 block) + new headers ---
$ws.Range("A43").Value = $ws.Range("A33").Value2
$ws.Range("B43").Value = -0.20082
$ws.Range("D43").Value = $ws.Range("D33").Value2
$ws.Range("E43").Value = $ws.Range("E33").Value2
$ws.Range("H43").Value = $ws.Range("H33").Value2
$ws.Range("I43").Value = $ws.Range("I33").Value2

# --- Row 44 ---
$ws.Range("A44").Value = $ws.Range("A34").Value2
$ws.Range("B44").Value = 7.001792
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 5.8
$ws.Range("E44").Value = -92
$ws.Range("F44").Formula = '=E44/D44'
$ws.Range("H44").Value = 3.1
$ws.Range("I44").Value = 118
$ws.Range("J44").Formula = '=I44/H44'

# --- Row 45 ---
$ws.Range("A45").Value = $ws.Range("A35").Value2
$ws.Range("B45").Value = 0.256173
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 6
$ws.Range("E45").Value = -124
$ws.Range("H45").Value = 3.1
$ws.Range("I45").Value = 122

# --- Row 46 ---
$ws.Range("A46").Value = $ws.Range("A36").Value2
$ws.Range("B46").Value = 0.91102
$ws.Range("C46").Value = 3
$ws.Range("D46").Value = 5.9
$ws.Range("E46").Value = -112
$ws.Range("H46").Value = 3.6
$ws.Range("I46").Value = 152

# --- Row 47 ---
$ws.Range("A47").Value = $ws.Range("A37").Value2
$ws.Range("B47").Value = 0.245392
$ws.Range("C47").Value = 4
$ws.Range("D47").Value = 5.9
$ws.Range("E47").Value = -93
$ws.Range("H47").Value = 3.7
$ws.Range("I47").Value = 135

# --- Row 48 ---
$ws.Range("A48").Value = $ws.Range("A38").Value2
$ws.Range("B48").Value = 0.348521
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 6
$ws.Range("E48").Value = -102
$ws.Range("H48").Value = 3.8
$ws.Range("I48").Value = 155

# Shared formulas for F45:F48 and J45:J48
$ws.Range("F45:F48").Formula = '=E45/D45'
$ws.Range("J45:J48").Formula = '=I45/H45'

# --- Row 49 ---
$ws.Range("A49").Value = $ws.Range("A39").Value2
$ws.Range("B49").Value = -0.036479999999999999
$ws.Range("F49").Formula = '=AVERAGE(F44:F48)'
$ws.Range("J49").Formula = '=AVERAGE(J44:J48)'

# --- Row 50 ---
$ws.Range("A50").Value = $ws.Range("A40").Value2
$ws.Range("B50").Value = 6.317713

# L column formulas (rows 43-46), non-shared individual formulas
$ws.Range("L43").Formula = '=B43*$F$49'
$ws.Range("L44").Formula = '=B44*$F$49'
$ws.Range("L45").Formula = '=B45*$F$49'
$ws.Range("L46").Formula = '=B46*$F$49'

# M column formulas (rows 47-50), non-shared individual formulas
$ws.Range("M47").Formula = '=B47*$J$49'
$ws.Range("M48").Formula = '=B48*$J$49'
$ws.Range("M49").Formula = '=B49*$J$49'
$ws.Range("M50").Formula = '=B50*$J$49'

# Final view state: scroll/select as in the edited workbook
$ws.Range("L43").Select() | Out-Null
